$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.849.80"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.084.54"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'233.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'59.30"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.396"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "'0.0791"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("D12").Value = "2.392.03"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'14.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "'0.775"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "'5.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "2.109.64"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "37.782.40"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'71.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").Value = "'228.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "'171.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'9.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.24%  "
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "'19.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'5.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").Value = "'99.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.63%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "1.450.53"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'1.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'4.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'7.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "2.277.32"
$ws.Range("E51").Value = "  -0.40%  "
